# OW-535: rename shared-string value "p1" -> "p1a" (cell AP2 on the
# IRS-Bilateral sheet) plus the incidental view/column-width drift that
# came along with the resave in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Content change -----------------------------------------------
# AP2 held the shared string "p1"; it becomes "p1a".
$ws.Range("AP2").Value = "p1a"

# --- 2. Selection / active cell ---------------------------------------
$ws.Range("AM14").Select()

# --- 3. Column width drift ---------------------------------------------
# Every column shifted slightly (a by-product of Excel recalculating
# "best fit" widths on save) and column C picked up an explicit width
# entry. ColumnWidth is expressed in characters; the values below are
# chosen so the persisted <col width="..."/> matches the target file.
$widths = @{
    1  = 8.333333333333332
    2  = 14.666666666666666
    3  = 7.666666666666667
    4  = 5.833333333333334
    6  = 10.833333333333332
    7  = 9.833333333333332
    8  = 9.5
    9  = 8.0
    10 = 5.666666666666667
    11 = 11.833333333333332
    12 = 9.5
    13 = 8.833333333333332
    14 = 14.333333333333332
    15 = 26.666666666666668
    16 = 16.5
    17 = 14.5
    18 = 10.166666666666666
    19 = 17.0
    20 = 16.166666666666668
    21 = 15.666666666666666
    22 = 14.333333333333332
    23 = 14.0
    24 = 15.333333333333334
    25 = 9.5
    26 = 8.833333333333332
    27 = 14.333333333333332
    28 = 26.666666666666668
    29 = 16.5
    30 = 14.5
    31 = 14.333333333333332
    32 = 17.0
    33 = 16.166666666666668
    34 = 15.666666666666666
    35 = 14.333333333333332
    36 = 14.0
    37 = 15.333333333333334
    38 = 14.5
    39 = 14.5
    40 = 11.333333333333332
    41 = 10.166666666666666
    42 = 8.833333333333332
}

foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col]
}
